$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "30.646.13"
$ws.Range('E2').Value = "  +0.51%  "
$ws.Range('D3').Value = "2.115.16"
$ws.Range('E3').Value = "  +0.19%  "
$ws.Range('E4').Value = "  +1.04%  "
$ws.Range('D5').Value = "'339.98"
$ws.Range('E5').Value = "  +1.52%  "
$ws.Range('E6').Value = "  +0.99%  "
$ws.Range('E7').Value = "  -0.03%  "
$ws.Range('D8').Value = "'0.4506"
$ws.Range('E8').Value = "  +0.40%  "
$ws.Range('E9').Value = "  +0.82%  "
$ws.Range('D10').Value = "'0.09082"
$ws.Range('E10').Value = "  +0.23%  "
$ws.Range('E11').Value = "  -0.12%  "
$ws.Range('D12').Value = "'24.37"
$ws.Range('E12').Value = "  -0.47%  "
$ws.Range('D13').Value = "2.125.75"
$ws.Range('E13').Value = "  +1.11%  "
$ws.Range('D14').Value = "'6.797"
$ws.Range('E14').Value = "  +0.12%  "
$ws.Range('D15').Value = "'8.069"
$ws.Range('E15').Value = "  +3.15%  "
$ws.Range('D16').Value = "'97.81"
$ws.Range('E16').Value = "  +0.99%  "
$ws.Range('D17').Value = "'0.00001163"
$ws.Range('E17').Value = "  +2.97%  "
$ws.Range('D18').Value = "'1.013"
$ws.Range('E18').Value = "  +1.02%  "
$ws.Range('D19').Value = "'0.06695"
$ws.Range('E19').Value = "  +1.08%  "
$ws.Range('D20').Value = "'19.37"
$ws.Range('E20').Value = "  -0.49%  "
$ws.Range('D21').Value = "'1.011"
$ws.Range('E21').Value = "  +0.98%  "
$ws.Range('E22').Value = "  +1.47%  "
$ws.Range('D23').Value = "30.731.12"
$ws.Range('E23').Value = "  +0.63%  "
$ws.Range('D24').Value = "'12.81"
$ws.Range('E24').Value = "  +3.27%  "
$ws.Range('D25').Value = "'2.377"
$ws.Range('E25').Value = "  +0.93%  "
$ws.Range('D26').Value = "2.367.49"
$ws.Range('E26').Value = "  +0.71%  "
$ws.Range('D27').Value = "'22.36"
$ws.Range('E27').Value = "  -0.26%  "
$ws.Range('D28').Value = "'165.26"
$ws.Range('E28').Value = "  +1.12%  "
$ws.Range('D29').Value = "'2.562"
$ws.Range('E29').Value = "  -1.29%  "
$ws.Range('D30').Value = "'135.73"
$ws.Range('E30').Value = "  +2.14%  "
$ws.Range('D31').Value = "'1.198"
$ws.Range('E31').Value = "  -0.23%  "
$ws.Range('E32').Value = "  -0.08%  "
$ws.Range('D33').Value = "'6.386"
$ws.Range('E33').Value = "  +3.50%  "
$ws.Range('D34').Value = "'1.634"
$ws.Range('E34').Value = "  -1.92%  "
$ws.Range('D35').Value = "'3.942"
$ws.Range('E35').Value = "  +0.17%  "
$ws.Range('E36').Value = "  -2.21%  "
$ws.Range('D37').Value = "'5.916"
$ws.Range('E37').Value = "  +6.10%  "
$ws.Range('D38').Value = "'0.02658"
$ws.Range('E38').Value = "  +2.88%  "
$ws.Range('D39').Value = "'0.06838"
$ws.Range('E39').Value = "  -0.04%  "
$ws.Range('D40').Value = "'0.2324"
$ws.Range('E40').Value = "  +0.63%  "
$ws.Range('D41').Value = "'12.62"
$ws.Range('E41').Value = "  -1.49%  "
$ws.Range('D42').Value = "'0.6882"
$ws.Range('E42').Value = "  -0.94%  "
$ws.Range('E43').Value = "  +1.22%  "
$ws.Range('D44').Value = "'15.03"
$ws.Range('E44').Value = "  +6.49%  "
$ws.Range('D45').Value = "'0.6433"
$ws.Range('E45').Value = "  +0.36%  "
$ws.Range('D46').Value = "'2.317"
$ws.Range('E46').Value = "  -2.92%  "
$ws.Range('E47').Value = "  +13.99%  "
$ws.Range('D48').Value = "'3.708"
$ws.Range('E48').Value = "  +1.05%  "
$ws.Range('E49').Value = "  +0.38%  "
$ws.Range('B50').Value = "Aave"
$ws.Range('C50').Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range('D50').Value = "'82.88"
$ws.Range('E50').Value = "  -0.69%  "
$ws.Range('B51').Value = "Cronos"
$ws.Range('C51').Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range('D51').Value = "'0.07313"
$ws.Range('E51').Value = "  +3.10%  "
